# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the associated handoff date/time stamps on the Overview, zh-cn and de-de
# sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-28-20 04:28:40"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("E2").Value = "2016-03-20 04:28:37"

# --- de-de sheet ---
$wsDe = $wb.Worksheets("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("E2").Value = "2016-03-20 04:28:40"
